$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.152.12"
$ws.Range("E2").Value = "  +6.95%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.570.53"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "416.66"
$ws.Range("E5").Value = "  +0.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.51"
$ws.Range("E6").Value = "  -0.34%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.652"
$ws.Range("E7").Value = "  +3.76%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.561.25"
$ws.Range("E8").Value = "  +2.72%  "
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.769"
$ws.Range("E10").Value = "  +5.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.174"
$ws.Range("E11").Value = "  +13.85%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000335"
$ws.Range("E12").Value = "  +47.93%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "42.39"
$ws.Range("E13").Value = "  -0.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.03"
$ws.Range("E14").Value = "  +2.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.137.31"
$ws.Range("E15").Value = "  +2.81%  "
$ws.Range("E16").Value = "  -0.18%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "20.44"
$ws.Range("E17").Value = "  -0.65%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.568.29"
$ws.Range("E18").Value = "  +2.98%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.14"
$ws.Range("E19").Value = "  +5.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "67.027.69"
$ws.Range("E20").Value = "  +6.88%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.34"
$ws.Range("E21").Value = "  -2.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "456.70"
$ws.Range("E22").Value = "  -1.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "88.05"
$ws.Range("E23").Value = "  -2.71%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.10"
$ws.Range("E24").Value = "  -5.76%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.44"
$ws.Range("E25").Value = "  +1.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.37"
$ws.Range("E26").Value = "  +1.78%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.14"
$ws.Range("E27").Value = "  -6.33%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "34.71"
$ws.Range("E28").Value = "  +4.13%  "
$ws.Range("E29").Value = "  +1.40%  "
$ws.Range("E30").Value = "  +4.45%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.37"
$ws.Range("E31").Value = "  +1.91%  "
$ws.Range("E32").Value = "  +4.81%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.38"
$ws.Range("E33").Value = "  -2.82%  "
$ws.Range("E34").Value = "  -4.36%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "41.21"
$ws.Range("E35").Value = "  +0.51%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "56.54"
$ws.Range("E37").Value = "  -2.81%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0495"
$ws.Range("E38").Value = "  +0.87%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0₃0726"
$ws.Range("E39").Value = "  +24.67%  "
$ws.Range("E40").Value = "  +8.86%  "
$ws.Range("E41").Value = "  -0.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.05"
$ws.Range("E42").Value = "  -0.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "149.05"
$ws.Range("E43").Value = "  -0.13%  "
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.27"
$ws.Range("E45").Value = "  -2.40%  "
$ws.Range("E46").Value = "  -3.69%  "
$ws.Range("E47").Value = "  -2.34%  "
$ws.Range("E48").Value = "  -3.96%  "
$ws.Range("E49").Value = "  -1.68%  "
$ws.Range("E50").Value = "  +14.92%  "
$ws.Range("E51").Value = "  -4.75%  "
